$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "2016-03-10 16:44:01"

$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "2016-03-10 16:44:06"
